# Budget.xlsx edit script
# Summary of the change (per commit message "JSON files are now automatically
# synced with Excel"): the Inventory sheet was updated with the store
# ("Location") each grocery item was bought at, a couple of previously blank
# rows (Zwiebel / Knoblauch) were filled in, the Dill quantity was corrected,
# and "Edeka " (trailing-space typo) was renamed to "Edeka". The active
# worksheet/selection also moved from the Recipes sheet to the Inventory sheet.

$wb = $excel.ActiveWorkbook

$inventory = $wb.Worksheets.Item("Inventory")
$recipes = $wb.Worksheets.Item("Recipes")

# --- Inventory sheet data edits -------------------------------------------

# Fix "Edeka " (trailing space) -> "Edeka" for Bread's store/location.
$inventory.Range("F4").Value = "Edeka"

# Fill in the "Location" column for rows that previously had none.
$inventory.Range("F5").Value = "Aldi"     # Walnut
$inventory.Range("F6").Value = "Aldi"     # Frischkase
$inventory.Range("F7").Value = "Edeka"    # Tea
$inventory.Range("F9").Value = "Aldi"     # Tomatenmark

# Zwiebel (row 10): add Minimum Purchase price / Category / Location.
$inventory.Range("D10").Value = 1.5
$inventory.Range("E10").Value = "Vegetables"
$inventory.Range("F10").Value = "Aldi"

# Knoblauch (row 11): previously almost empty, now fully filled in.
$inventory.Range("B11").Value = 20
$inventory.Range("C11").Value = 1
$inventory.Range("D11").Value = 1.5
$inventory.Range("E11").Value = "Vegetables"
$inventory.Range("F11").Value = "Aldi"

# Dill (row 12): weight corrected from 15 to 13, and location added.
$inventory.Range("B12").Value = 13
$inventory.Range("F12").Value = "Edeka"

# Column D (Purchase Price) widened to fit the new longer currency values.
$inventory.Columns.Item(4).ColumnWidth = 13.711495535714286

# --- Selection / active sheet changes --------------------------------------

# Recipes was the active sheet before; its saved selection moves from C2 to B3
# and it is no longer the tab shown when the workbook is reopened.
$recipes.Activate()
$recipes.Range("B3").Select()

# Inventory becomes the active sheet, with the cursor on A3 (it had been on
# D10 before).
$inventory.Activate()
$inventory.Range("A3").Select()
